$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 307, shifting existing rows 307:340 down to 308:341
$ws.Rows.Item(307).Insert()

# Populate the newly inserted row 307 with the new weekly record.
# Columns A,B,C,E,F,G,H,I,J,K,Q,T mirror the constant values already used
# throughout this Frambuesa / Mercado Mayorista Lo Valledor block.
$ws.Cells.Item(307, 1).Value = 6
$ws.Cells.Item(307, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(307, 3).Value = "Metropolitana"
$ws.Cells.Item(307, 4).Value = (Get-Date -Year 2023 -Month 3 -Day 28 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(307, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(307, 5).Value = 13
$ws.Cells.Item(307, 6).Value = "Fruta"
$ws.Cells.Item(307, 7).Value = 100101
$ws.Cells.Item(307, 8).Value = "Berries"
$ws.Cells.Item(307, 9).Value = 100101004
$ws.Cells.Item(307, 10).Value = "Frambuesa"
$ws.Cells.Item(307, 11).Value = "Sin especificar"
$ws.Cells.Item(307, 12).Value = "Primera"
$ws.Cells.Item(307, 13).Value = 250
$ws.Cells.Item(307, 14).Value = 7000
$ws.Cells.Item(307, 15).Value = 7000
$ws.Cells.Item(307, 16).Value = 7000
$ws.Cells.Item(307, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(307, 18).Value = "Región del Maule"
$ws.Cells.Item(307, 19).Value = 3500
$ws.Cells.Item(307, 20).Value = 2
